$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.896.97'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.116.31'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '347.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5203'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4453'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.36'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09339'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  +0.93%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.450'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("D14").Value = '2.104.30'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.867'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.20'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.18%  '
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.007'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06672'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.295'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").Value = '29.931.26'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.329'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '2.349.38'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.554'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.792'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1056'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.233'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.969'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.471'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02600'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06856'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7024'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2250'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.335'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6815'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.349'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.09%  '
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.636'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000356'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.242'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.57%  '
$ws.Range("E51").Value = '  +0.12%  '
